# Apply the update described by the commit: refresh the "Date" metadata
# value and add properly-capitalized Display values for three of the
# "Missing - ..." concepts on the Concepts sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value (row 8, column B) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# --- Concepts sheet: capitalize Display text for three rows ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C3").Value = "Missing - Restricted Access"
$concepts.Range("C4").Value = "Missing - Not Provided"
$concepts.Range("C5").Value = "Missing - Not Collected"
